# Added connection between cost and demand
# Added variable and fixed costs to connect to demand model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Q_toBoil (B3) with new computed value; T_boil (B4) recalculates via
# its existing formula =B3/B2.
$ws.Range("B3").Value = 212831.54173463001

$wb.Application.Calculate()
